# Apply the update described by the diff:
#  1. Change column C ("Forandrad" / changed date) from 45182 to 45184 for all
#     data rows (rows 2..171).
#  2. Give row 171 an explicit row height (15, custom height) so it matches
#     the formatting of the other data rows.
#  3. Append a brand-new row 172 with the new cutting notification record
#     "A 42937-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = 45182
$newDate = 45184

# --- 1. Update column C for every existing data row (rows 2..171) ----------
for ($r = 2; $r -le 171; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value2 = $newDate
    }
}

# --- 2. Make row 171 carry an explicit custom row height --------------------
$ws.Rows(171).RowHeight = 15

# --- 3. Append the new row 172 ----------------------------------------------
$newRow = 172

$ws.Cells.Item($newRow, 1).Value = "A 42937-2023"

$ws.Cells.Item($newRow, 2).Value2 = 45182
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item(171, 2).NumberFormat

$ws.Cells.Item($newRow, 3).Value2 = $newDate
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item(171, 3).NumberFormat

$ws.Cells.Item($newRow, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item($newRow, 5).Value = "HÖÖR"

$ws.Cells.Item($newRow, 7).Value2 = 2.1
$ws.Cells.Item($newRow, 8).Value2 = 0
$ws.Cells.Item($newRow, 9).Value2 = 0
$ws.Cells.Item($newRow, 10).Value2 = 0
$ws.Cells.Item($newRow, 11).Value2 = 0
$ws.Cells.Item($newRow, 12).Value2 = 0
$ws.Cells.Item($newRow, 13).Value2 = 0
$ws.Cells.Item($newRow, 14).Value2 = 0
$ws.Cells.Item($newRow, 15).Value2 = 0
$ws.Cells.Item($newRow, 16).Value2 = 0
$ws.Cells.Item($newRow, 17).Value2 = 0

$ws.Cells.Item($newRow, 18).WrapText = $ws.Cells.Item(171, 18).WrapText
$ws.Cells.Item($newRow, 18).Value = ""
